$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as text so values like "1.000" or "30.493.96"
# are not auto-converted/reformatted as numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '30.493.96'
$ws.Range("E2").Value = '  +1.03%  '

# Row 3
$ws.Range("D3").Value = '1.879.57'
$ws.Range("E3").Value = '  +0.91%  '

# Row 4
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  +0.20%  '

# Row 5
$ws.Range("D5").Value = '246.29'
$ws.Range("E5").Value = '  +5.19%  '

# Row 6
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  +0.16%  '

# Row 7
$ws.Range("D7").Value = '0.4752'
$ws.Range("E7").Value = '  +1.55%  '

# Row 8
$ws.Range("D8").Value = '0.2898'
$ws.Range("E8").Value = '  +2.05%  '

# Row 9
$ws.Range("D9").Value = '0.06514'
$ws.Range("E9").Value = '  +0.49%  '

# Row 10
$ws.Range("D10").Value = '21.65'
$ws.Range("E10").Value = '  +2.07%  '

# Row 11
$ws.Range("D11").Value = '0.07734'
$ws.Range("E11").Value = '  -0.17%  '

# Row 12
$ws.Range("B12").Value = 'Polygon'
$ws.Range("C12").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D12").Value = '0.7416'
$ws.Range("E12").Value = '  +8.60%  '

# Row 13
$ws.Range("B13").Value = 'Litecoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D13").Value = '96.99'
$ws.Range("E13").Value = '  +3.21%  '

# Row 14
$ws.Range("D14").Value = '1.880.55'
$ws.Range("E14").Value = '  +0.96%  '

# Row 15
$ws.Range("D15").Value = '5.117'
$ws.Range("E15").Value = '  +1.40%  '

# Row 16
$ws.Range("D16").Value = '273.47'
$ws.Range("E16").Value = '  +0.52%  '

# Row 17
$ws.Range("D17").Value = '30.495.90'
$ws.Range("E17").Value = '  +1.14%  '

# Row 18
$ws.Range("D18").Value = '13.66'
$ws.Range("E18").Value = '  +2.22%  '

# Row 19
$ws.Range("D19").Value = '0.000007563'
$ws.Range("E19").Value = '  -0.19%  '

# Row 20
$ws.Range("D20").Value = '1.000'
$ws.Range("E20").Value = '  +0.08%  '

# Row 21
$ws.Range("D21").Value = '2.131.51'
$ws.Range("E21").Value = '  +0.22%  '

# Row 22
$ws.Range("D22").Value = '1.001'
$ws.Range("E22").Value = '  +0.28%  '

# Row 23
$ws.Range("D23").Value = '5.260'
$ws.Range("E23").Value = '  +2.26%  '

# Row 24
$ws.Range("D24").Value = '6.169'
$ws.Range("E24").Value = '  +1.05%  '

# Row 25
$ws.Range("E25").Value = '  -1.06%  '

# Row 26
$ws.Range("D26").Value = '164.53'
$ws.Range("E26").Value = '  -0.69%  '

# Row 27
$ws.Range("D27").Value = '18.89'
$ws.Range("E27").Value = '  +1.93%  '

# Row 28
$ws.Range("D28").Value = '1.951'
$ws.Range("E28").Value = '  +3.02%  '

# Row 29
$ws.Range("D29").Value = '1.374'
$ws.Range("E29").Value = '  +1.02%  '

# Row 30
$ws.Range("D30").Value = '0.09978'
$ws.Range("E30").Value = '  +1.25%  '

# Row 31
$ws.Range("D31").Value = '1.512'
$ws.Range("E31").Value = '  +3.94%  '

# Row 32
$ws.Range("D32").Value = '4.320'
$ws.Range("E32").Value = '  +1.91%  '

# Row 33
$ws.Range("D33").Value = '4.064'
$ws.Range("E33").Value = '  +1.84%  '

# Row 34
$ws.Range("D34").Value = '0.04765'
$ws.Range("E34").Value = '  +2.36%  '

# Row 35
$ws.Range("D35").Value = '1.124'
$ws.Range("E35").Value = '  +0.51%  '

# Row 36
$ws.Range("D36").Value = '0.6981'
$ws.Range("E36").Value = '  +1.13%  '

# Row 37
$ws.Range("D37").Value = '2.716'
$ws.Range("E37").Value = '  +0.49%  '

# Row 38
$ws.Range("D38").Value = '0.01870'
$ws.Range("E38").Value = '  +1.84%  '

# Row 39
$ws.Range("D39").Value = '2.737'
$ws.Range("E39").Value = '  -0.26%  '

# Row 40
$ws.Range("D40").Value = '6.340'
$ws.Range("E40").Value = '  +0.16%  '

# Row 41
$ws.Range("D41").Value = '70.21'
$ws.Range("E41").Value = '  -1.24%  '

# Row 42
$ws.Range("D42").Value = '1.919'
$ws.Range("E42").Value = '  +1.85%  '

# Row 43
$ws.Range("D43").Value = '0.4175'
$ws.Range("E43").Value = '  +2.77%  '

# Row 44
$ws.Range("D44").Value = '0.9999'
$ws.Range("E44").Value = '  +0.14%  '

# Row 45
$ws.Range("D45").Value = '0.8390'
$ws.Range("E45").Value = '  +1.21%  '

# Row 46
$ws.Range("D46").Value = '102.55'
$ws.Range("E46").Value = '  +0.07%  '

# Row 47
$ws.Range("D47").Value = '9.316'
$ws.Range("E47").Value = '  +3.69%  '

# Row 48
$ws.Range("D48").Value = '7.086'
$ws.Range("E48").Value = '  +1.70%  '

# Row 49
$ws.Range("D49").Value = '35.42'
$ws.Range("E49").Value = '  +4.56%  '

# Row 50
$ws.Range("D50").Value = '928.30'
$ws.Range("E50").Value = '  -0.54%  '

# Row 51
$ws.Range("D51").Value = '0.05621'
$ws.Range("E51").Value = '  +0.88%  '
